$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.003.01'
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("D3").Value = '2.924.45'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.51'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.63%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.92'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.53%  '
$ws.Range("E10").Value = '  -0.24%  '
$ws.Range("E11").Value = '  -1.11%  '
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.59'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.56%  '
$ws.Range("E14").Value = '  +0.18%  '
$ws.Range("D15").Value = '3.408.74'
$ws.Range("E15").Value = '  +0.56%  '
$ws.Range("D16").Value = '60.919.96'
$ws.Range("E16").Value = '  +0.39%  '
$ws.Range("D18").Value = '2.925.24'
$ws.Range("E18").Value = '  +0.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '432.26'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.37'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.61%  '
$ws.Range("E21").Value = '  -0.39%  '
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '81.45'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("E25").Value = '  -0.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.87'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("E28").Value = '  +5.37%  '
$ws.Range("E29").Value = '  +0.29%  '
$ws.Range("E30").Value = '  -2.08%  '
$ws.Range("E31").Value = '  +0.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.108'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.27%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").Value = '0.0₃0856'
$ws.Range("E34").Value = '  -0.70%  '
$ws.Range("E35").Value = '  +0.69%  '
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.99'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.29%  '
$ws.Range("E38").Value = '  -0.72%  '
$ws.Range("E39").Value = '  -3.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.57'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.77%  '
$ws.Range("E41").Value = '  -3.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.03'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '380.83'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.58%  '
$ws.Range("D44").Value = '2.700.84'
$ws.Range("E44").Value = '  +1.26%  '
$ws.Range("E45").Value = '  -1.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '133.97'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.84%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.82'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.95%  '
$ws.Range("E49").Value = '  -0.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.01'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.54%  '
$ws.Range("E51").Value = '  -0.60%  '
